$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.142531633377075
$ws.Range("B1").Value = 1.254926204681396
$ws.Range("C1").Value = 1.488237619400024
$ws.Range("D1").Value = 2.629266738891602
$ws.Range("E1").Value = 4.470052242279053
